$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting/values of the last existing row (22) into the new
# row (23), then overwrite just the two cells that actually carry new data
# (Task No and Task). This keeps the date-like "End date" text cell (F) and
# the numeric "Start date" (E) exactly as they were on row 22, including
# their cell styles, instead of letting Excel reinterpret a typed string as
# a date value.
$ws.Range("A22:G22").Copy($ws.Range("A23:G23"))
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Get DSA Code of XMED from Niha's PC vah vah"

# Grow the table ("Table2") so the new row is included, keeping the
# AutoFilter range in sync.
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A1:H23"))

# Match Excel's behaviour of moving the active selection to the newly
# entered cell.
$ws.Range("B23").Select() | Out-Null
